# Applies the numeric corrections to the "alpha diversity" regression
# table (Tables/alphaTable2.docx) as described by the commit:
#   "Dove into alpha diversity to remember about how removing spikes
#    leads to lower sample reads"
#
# The table has columns: Metric | Term | Estimate | Standard Error | T Value | p
# Scientific-notation cells ("m x 10^e") are stored as two runs: the
# mantissa ("m x 10") in a normal run, and the exponent (e) in a
# superscript run. We edit each run precisely so we never disturb the
# other run's character formatting.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nbsp = [string][char]0x00A0

# Replaces the text of a whole table cell. Only safe for cells that
# contain a single run (e.g. plain "T Value"/"p" numbers) because it
# sets .Text on the cell's full Range (which Word re-terminates with
# the cell mark automatically).
function Set-SingleRunCell($row, $col, $newText) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}

# Replaces just the leading "mantissa" run of a scientific-notation
# cell (e.g. "1.5 x 10" in "1.5 x 10^3"), leaving the trailing
# superscript exponent run untouched.
function Set-MantissaRun($row, $col, $oldMantissa, $newMantissa) {
    $cell = $t.Cell($row, $col)
    $start = $cell.Range.Start
    $sub = $d.Range($start, $start + $oldMantissa.Length)
    if ($sub.Text -ne $oldMantissa) {
        throw "Mantissa mismatch at row $row col $col : expected '$oldMantissa' got '$($sub.Text)'"
    }
    $sub.Text = $newMantissa
}

# Replaces just the trailing superscript exponent run of a
# scientific-notation cell, leaving the mantissa run untouched.
# mantissaLen must be the length (in characters) of the mantissa run
# *before* any edit made to it in the same cell.
function Set-ExponentRun($row, $col, $mantissaLen, $oldExponent, $newExponent) {
    $cell = $t.Cell($row, $col)
    $start = $cell.Range.Start
    $sub = $d.Range($start + $mantissaLen, $start + $mantissaLen + $oldExponent.Length)
    if ($sub.Text -ne $oldExponent) {
        throw "Exponent mismatch at row $row col $col : expected '$oldExponent' got '$($sub.Text)'"
    }
    $sub.Text = $newExponent
}

# ---- Diversity (Shannon H) ----------------------------------------
# Intercept
Set-MantissaRun 14 3 ("1.5" + $nbsp + "x" + $nbsp + "10") ("1.4" + $nbsp + "x" + $nbsp + "10")
Set-SingleRunCell 14 5 "1.84"
Set-SingleRunCell 14 6 "0.071"

# log(Size Class)
Set-MantissaRun 15 3 ("2.0" + $nbsp + "x" + $nbsp + "10") ("1.9" + $nbsp + "x" + $nbsp + "10")
Set-SingleRunCell 15 5 "4.07"

# log(Size Class)^2
Set-MantissaRun 16 3 ("-3.8" + $nbsp + "x" + $nbsp + "10") ("-3.5" + $nbsp + "x" + $nbsp + "10")
Set-MantissaRun 16 4 ("8.8" + $nbsp + "x" + $nbsp + "10") ("8.7" + $nbsp + "x" + $nbsp + "10")
Set-SingleRunCell 16 5 "-4.00"

# Latitude
Set-MantissaRun 17 3 ("-7.9" + $nbsp + "x" + $nbsp + "10") ("-7.4" + $nbsp + "x" + $nbsp + "10")
Set-MantissaRun 17 4 ("4.1" + $nbsp + "x" + $nbsp + "10") ("4.0" + $nbsp + "x" + $nbsp + "10")
Set-SingleRunCell 17 5 "-1.83"
Set-SingleRunCell 17 6 "0.072"

# Latitude^2 (mantissa AND exponent change: 1.0 x 10^0 -> 9.6 x 10^-1)
$mantissaOld = "1.0" + $nbsp + "x" + $nbsp + "10"
$mantissaNew = "9.6" + $nbsp + "x" + $nbsp + "10"
Set-ExponentRun 18 3 $mantissaOld.Length "0" "-1"
Set-MantissaRun 18 3 $mantissaOld $mantissaNew
Set-SingleRunCell 18 5 "1.83"
Set-SingleRunCell 18 6 "0.072"

# Depth
Set-MantissaRun 19 3 ("1.6" + $nbsp + "x" + $nbsp + "10") ("1.5" + $nbsp + "x" + $nbsp + "10")
Set-SingleRunCell 19 5 "0.78"
Set-SingleRunCell 19 6 "0.436"

# ---- Evenness (Pielou J) ------------------------------------------
# Intercept
Set-SingleRunCell 20 5 "-0.72"
Set-SingleRunCell 20 6 "0.476"

# log(Size Class)
Set-SingleRunCell 21 6 "0.052"

# log(Size Class)^2
Set-MantissaRun 22 3 ("5.7" + $nbsp + "x" + $nbsp + "10") ("5.8" + $nbsp + "x" + $nbsp + "10")
Set-SingleRunCell 22 6 "0.068"

# Latitude
Set-SingleRunCell 23 6 "0.477"

# Latitude^2
Set-SingleRunCell 24 6 "0.478"

# Depth
Set-SingleRunCell 25 6 "0.718"
